$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.998.49"
$ws.Range("E2").Value = "  -2.70%  "
$ws.Range("D3").Value = "1.794.61"
$ws.Range("E3").Value = "  -3.12%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.59%  "
$ws.Range("D7").Value = "'0.4175"
$ws.Range("E7").Value = "  -2.69%  "
$ws.Range("D8").Value = "'0.3549"
$ws.Range("E8").Value = "  -4.47%  "
$ws.Range("D9").Value = "'0.07071"
$ws.Range("E9").Value = "  -4.15%  "
$ws.Range("D10").Value = "'0.8421"
$ws.Range("E10").Value = "  -4.06%  "
$ws.Range("D11").Value = "'20.11"
$ws.Range("E11").Value = "  -4.45%  "
$ws.Range("D12").Value = "1.839.42"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").Value = "'5.271"
$ws.Range("E13").Value = "  -3.45%  "
$ws.Range("D14").Value = "'6.332"
$ws.Range("E14").Value = "  -4.00%  "
$ws.Range("D15").Value = "'0.06730"
$ws.Range("E15").Value = "  -3.06%  "
$ws.Range("D16").Value = "'1.009"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").Value = "'79.72"
$ws.Range("E17").Value = "  -1.73%  "
$ws.Range("D18").Value = "'0.000008700"
$ws.Range("E18").Value = "  -4.46%  "
$ws.Range("D19").Value = "'1.005"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").Value = "'14.97"
$ws.Range("E20").Value = "  -3.73%  "
$ws.Range("D21").Value = "27.277.10"
$ws.Range("E21").Value = "  -1.59%  "
$ws.Range("D22").Value = "'5.046"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("D23").Value = "'10.91"
$ws.Range("E23").Value = "  -1.37%  "
$ws.Range("D24").Value = "2.028.39"
$ws.Range("E24").Value = "  -3.75%  "
$ws.Range("D25").Value = "'1.940"
$ws.Range("E25").Value = "  -1.51%  "
$ws.Range("D26").Value = "'152.94"
$ws.Range("E26").Value = "  -1.20%  "
$ws.Range("D28").Value = "'4.969"
$ws.Range("E28").Value = "  -6.89%  "
$ws.Range("D29").Value = "'113.18"
$ws.Range("E29").Value = "  -2.06%  "
$ws.Range("D31").Value = "'0.08898"
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("D32").Value = "'0.7143"
$ws.Range("E32").Value = "  -8.92%  "
$ws.Range("D33").Value = "'2.857"
$ws.Range("E33").Value = "  -3.89%  "
$ws.Range("D34").Value = "'4.305"
$ws.Range("E34").Value = "  -6.57%  "
$ws.Range("D36").Value = "'1.070"
$ws.Range("E36").Value = "  -7.75%  "
$ws.Range("D37").Value = "'1.074"
$ws.Range("E37").Value = "  -3.26%  "
$ws.Range("D38").Value = "'0.01893"
$ws.Range("E38").Value = "  -3.45%  "
$ws.Range("D39").Value = "'0.05105"
$ws.Range("E39").Value = "  -5.94%  "
$ws.Range("D40").Value = "'0.1615"
$ws.Range("E40").Value = "  -3.94%  "
$ws.Range("D41").Value = "'0.4927"
$ws.Range("E41").Value = "  -5.39%  "
$ws.Range("D42").Value = "'2.579"
$ws.Range("E42").Value = "  -9.13%  "
$ws.Range("D43").Value = "'6.066"
$ws.Range("E43").Value = "  -10.44%  "
$ws.Range("D44").Value = "'7.984"
$ws.Range("E44").Value = "  -7.81%  "
$ws.Range("D45").Value = "'104.43"
$ws.Range("E45").Value = "  -2.90%  "
$ws.Range("D47").Value = "'10.20"
$ws.Range("E47").Value = "  -4.06%  "
$ws.Range("D48").Value = "'0.06302"
$ws.Range("E48").Value = "  -4.07%  "
$ws.Range("D49").Value = "'0.4477"
$ws.Range("E49").Value = "  -6.40%  "
$ws.Range("D50").Value = "'1.588"
$ws.Range("E50").Value = "  -4.53%  "
$ws.Range("D51").Value = "'61.80"
$ws.Range("E51").Value = "  -4.96%  "
$ws.Range("E27").Value = "  -2.69%  "
$ws.Range("E30").Value = "  -12.06%  "
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("E46").Value = "  +0.46%  "
$ws.Range("B5").Value = "USDC"
$ws.Range("C5").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D5").Value = "'1.005"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'307.59"
$ws.Range("E6").Value = "  -2.33%  "
